$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F4: mark bug #2 as CORREGIDO (copy formatting+value from an existing
#         "CORREGIDO" styled cell so the fill/font/border match exactly) ---
$ws.Range("F13").Copy($ws.Range("F4"))

# --- Row 18: new bug entry (#16) ---
# Fill in order C, B, E so new shared-string indices land in the same
# sequence as the target workbook.
$ws.Range("C18").Value() = "Método obtenerUltimoTorneoDelUsurio en DAOTorneo: CAMBIAR!!!!"
$ws.Range("B18").Value() = "Último Torneo del Usuario"
$ws.Range("E18").Value() = "login.aspx"

$ws.Range("A18").Value() = 16
$ws.Range("D18").Value() = "Facu"

# F18 gets the PENDIENTE text with a fresh style (copy PENDIENTE formatting
# from F17, then drop the vertical-centering so a new cellXf is minted).
$ws.Range("F17").Copy($ws.Range("F18"))
$ws.Range("F18").VerticalAlignment() = -4107

# C18 wraps its (long) description text.
$ws.Range("C18").WrapText() = $true

$ws.Rows.Item(18).RowHeight() = 30

# --- Selection / scroll position ---
$ws.Range("I18").Select()

Write-Host "done"
